$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION2")

# Update the evaluation level for indicator rows 14 and 15:
# change from "Logrado" to "Completamente logrado".
$ws.Range("C14").Value = "Completamente logrado"
$ws.Range("C15").Value = "Completamente logrado"

# Update the active selection on the sheet (was C24, now B24).
$ws.Activate()
$ws.Range("B24").Select()
